# "Completed array and graph"
# Adds a new worksheet "python PQ" (between "python DS" and "Login") that contains
# the existing Code/expectedOutcome sample rows plus four new array-algorithm
# practice problems (search, findMaxConsecutiveOnes, findNumbers, sortedSquares)
# together with their expected console output.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new sheet in the right slot (right after "python DS").
# ---------------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("Login")
$ws = $wb.Worksheets.Add($loginSheet)
$ws.Name = "python PQ"

# ---------------------------------------------------------------------------
# 2. Column sizing (wide column for the code snippets, narrower for results).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 212.08984375
$ws.Columns.Item(2).ColumnWidth = 35

# ---------------------------------------------------------------------------
# 3. Helper for building multi-line python snippets (PowerShell backtick-b is
#    a literal backspace, matching the blank "dedent" lines from the source).
# ---------------------------------------------------------------------------
$bs = [char]8

$searchCode = @(
  'def search(input_list, num):',
  'if(num in input_list):',
  'print("Element Found")',
  $bs,
  $bs,
  'else:',
  'print("Not Found")',
  $bs,
  $bs,
  $bs,
  $bs,
  'search([12, 23, 45, 67, 6, 90] , 12)'
) -join "`n"

$maxConsecutiveCode = @(
  'def findMaxConsecutiveOnes(nums) :',
  'count = 0',
  'result = 0',
  'for i in range(0, len(nums)):',
  'if (nums[i] == 0):',
  'count = 0',
  $bs,
  $bs,
  'else:',
  'count+= 1',
  $bs,
  $bs,
  'result = max(result, count)',
  $bs,
  $bs,
  'print(result)',
  $bs,
  $bs,
  'findMaxConsecutiveOnes([1,0,1,1,0,1])'
) -join "`n"

$findNumbersCode = @(
  'def findNumbers(nums):',
  'c=0',
  'for i in nums:',
  'j=str(i)',
  'x=len(j)',
  'if x%2==0:',
  'c=c+1',
  $bs,
  $bs,
  $bs,
  $bs,
  'print c',
  'return c',
  'findNumbers([12,345,2,6,7896])'
) -join "`n"

$sortedSquaresCode = @(
  'def sortedSquares(nums):',
  'squares_list = []',
  'for i in range(0, len(nums)):',
  'square = nums[i] * nums[i];',
  'squares_list.append(square)',
  $bs,
  $bs,
  'sorted_squares_list = sorted(squares_list)',
  'print sorted_squares_list;',
  'return sorted_squares_list;',
  'sortedSquares([-7,-3,2,3,11])'
) -join "`n"

# ---------------------------------------------------------------------------
# 4. Content. Rows 1-3 are the same Code/expectedOutcome rows already present
#    on the "LinkedList" sheet; rows 4-7 are the new practice questions.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "expectedOutcome"

$ws.Range("A2").Value = 'print("hello)abc'
$ws.Range("B2").Value = "SyntaxError: bad input on line 1"

$ws.Range("A3").Value = 'print("hello")'
$ws.Range("B3").Value = "hello"

$ws.Range("A4").Value = $searchCode
$ws.Range("B4").Value = "Element Found"

$ws.Range("A5").Value = $maxConsecutiveCode
$ws.Range("B5").Value = 2

$ws.Range("A6").Value = $findNumbersCode
$ws.Range("B6").Value = 2

$ws.Range("A7").Value = $sortedSquaresCode
$ws.Range("B7").Value = "[4, 9, 9, 49, 121]"

# Row 8 stays blank (matches the trailing blank row left behind in the source file)
# but still carries formatting, so touch it lightly to keep it in the sheet's
# used range (dimension A1:B8, like the source).
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = ""

# ---------------------------------------------------------------------------
# 5. Formatting: thin border around the "table" area, and a highlighted
#    left-aligned style for the new code cells (A4:A7).
# ---------------------------------------------------------------------------
$borderRange = $ws.Range("A1:A8")
$borderRange.Borders.LineStyle = 1
$borderRange.Borders.Weight = 2
$ws.Range("B1:B6").Borders.LineStyle = 1
$ws.Range("B1:B6").Borders.Weight = 2

$codeRange = $ws.Range("A4:A7")
$codeRange.Interior.Color = 16777215
$codeRange.Font.Color = 0
$codeRange.HorizontalAlignment = -4131

$ws.Range("B7").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 6. Selections, to mirror the saved workbook's view state.
# ---------------------------------------------------------------------------
$dsSheet = $wb.Worksheets.Item("python DS")
[void]$dsSheet.Activate()
[void]$dsSheet.Range("A3").Select()

$linkedListSheet = $wb.Worksheets.Item("LinkedList")
[void]$linkedListSheet.Activate()
[void]$linkedListSheet.Range("B3").Select()

$registerSheet = $wb.Worksheets.Item("Register")
[void]$registerSheet.Activate()
[void]$registerSheet.Range("A1:C9").Select()

[void]$ws.Activate()
[void]$ws.Range("B8").Select()

Write-Output "python PQ sheet added"
